$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing column D (Fact) values ---
$dChanges = @{
    13 = "3150000.0"
    14 = "2060.5333333333333"
    15 = "900000.0"
    17 = "9534.375"
    18 = "7560.0"
    20 = "8400.671999999999"
    22 = "8050000.0"
    23 = "2650000.0"
    24 = "1056000.0"
    25 = "1260000.0"
    26 = "303708.0"
    27 = "1194000.0"
    28 = "1674996.0"
    29 = "999996.0"
    30 = "1311000.0"
    31 = "nan"
    32 = "656454.0"
    33 = "1110000.0"
    34 = "1500000.0"
    35 = "774000.0"
    38 = "7845967.451"
    39 = "3652780.0"
    41 = "4125.0"
    44 = "1365.44"
    45 = "2240.0"
    46 = "1960.0"
    47 = "2000.6"
    48 = "3150.0"
    49 = "nan"
    50 = "nan"
    51 = "nan"
}

foreach ($row in $dChanges.Keys) {
    $cell = $ws.Cells.Item([int]$row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dChanges[$row]
    $cell.ClearFormats()
}

# --- New rows 52-64 appended at the bottom of the table ---
$newRows = @(
    @("Блоггеры", "None", "nan", "None", "None"),
    @("Максима Телеком ( Qvant)wi-fi.ru", "None", "300.0", "None", "None"),
    @("Блоггеры", "None", "nan", "None", "None"),
    @("Максима Телеком ( Qvant)wi-fi.ru", "None", "300.0", "None", "None"),
    @("Блоггеры", "None", "nan", "None", "None"),
    @("Максима Телеком ( Qvant)wi-fi.ru", "None", "300.0", "None", "None"),
    @("Smart TVGPMD", "None", "330.0", "None", "None"),
    @("Smart TVGPMD", "None", "330.0", "None", "None"),
    @("Smart TVИМХО", "None", "845000.0", "None", "None"),
    @("Блоггеры", "None", "nan", "None", "None"),
    @("Максима Телеком ( Qvant)wi-fi.ru", "None", "300.0", "None", "None"),
    @("Пакет XL Flex Rambler&Сo Desktop+Mobile Reach Video PMP", "None", "1350.0", "None", "None"),
    @("Пакет XL Flex Rambler&Сo Desktop+Mobile Reach Video PMP", "None", "1350.0", "None", "None")
)

$startRow = 52
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c]
        $cell.ClearFormats()
    }
}

$ws.Range("A1").Select()
